$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet from "My Series" to "Data"
$ws.Name = "Data"

# 2) Update the cached add-in metadata blob stored in the A1 cell comment
#    (gzip-compressed XML, base64-encoded, prefixed with the decompressed length)
$commentText = "Jx0AAB+LCAAAAAAAAAOlWVtvI0kV/istP4GE3d12ZjaJKr3yLVkLO45sh0z2BZW7K3GRdrfpqk7it0UCLVoWIYRm0XIVT4uQGEawKy0zXP7LapIZnvgLnLr0zXZ20mEUTbrO+U7VqVPnVhX07vXcNy5JxGgY7FXsmlUxSOCGHg3O9yoxP6vajyvvOqh77RL/CEd4TjiADZAK2O41o3uVGeeLXdO8urqqXTVqYXRu1i3LNp8M+mN3Rua4SgPGceCSSirlvV2q4qC2Nx8Qjj3MsZLcq/TGvVqbULcDtAEO8DmJaq2Y0YAw1g045ZQwIRkRzEm7M/iO2phTrz2u2chco2fIVkx9T+EKSEXXOFiWTOicOHXL3q5a29WGNbHru/bWrrVTs+r19xPBFIj6mPExiS6pKwljjucLKW5tNyy7bm9ZO8jcCIK5MgM4aOh7I3JJGfHaxPdZKYuY+gCbLoddlzOmhcycrJ7o4SocRHgxm1Duk3JqjAYtYx5oXbJJHLQfRsQF+z1IpUNyNYy0WSeLPnAnMxrxZQcvS891zEg0XAgjlRN1UCcMeNMnET9ewFkTD1wBGA6PYoLMO5iZUIcyF75pEBPPOcM+ywsVmOgkjC7YArvkEOLYFHNcBX6IPXA4ThmnbrboGgMdReECZoTFW6Hv7cOsGryBkc7cC8DEYtlWGF5k2m1iInmq8nzhTOeYJ/A1OhrPwqth4C/H8ZS5EZ0Sr9NK0Bt5SASklm7HjIdz0CIjIUXLUZbwDwJwlYw6xKVz7B/5YETmNGCWAgE1Yx6eUd4O/XgesESnFSo6gR1NyHW6w3SMhnC4gTB6GPSCBK/MvJFVFBiFV+ma6wxphBy5ydzkuNcZq+AO0JLjW+fIExG73Kc+1If8WeSoRa8YzwjhG11CcZBIhfui4jit5WE8n0J4TSHGLuWqDJkZH4Gfgq+DXo4FVaQqfyaWtSt/QI+UjbqBdzcuYSJYLreWYwNvhYRgT37Lx8EFUE8onx02k71s4CBlgTvx6zwEkbvw8VKSUyvlaagXuH7sEZUQesGZdFGhmzrUO9lojdSHGHcQDpaT5QLyMqO7HD72KlCpdxmPoBeoOG4YBzxaisyBTA19mwyLp4FcAPv3ljmLyPdjaEGW+3HgtkPv/qt5yjrHAeX31zCMI5UO7y8irScyY8w6ROQYmfTvLe+W2ROLSsHnAZmHAXXvb20wstDee8BGWBJV95YgKr7ujfehrquyJ2L93mIR9I9Q6Eot02QsdKl0Vh0eXk7evCNkOuQMxz70bhxK7Hmae1fJqMkuVjF5EjqO/CQDOqIzZtAau9685kLzINq/mhvOBcGEjvRkjMw8XnRALukG530cnMfQY6R5ZZWe5l9RHycRDpjYTtpSrKTizSCU5CnV6jgqeQ1j6QgqeYXAReYKDk3IfBFG2B+AYei+djvdL0ErMsB8pkdQ23ziJkY2M9FUqqhZovjbYLJIqW2IgNdpcoUoQWIvqgnPMBkNiV0OICz9NvbpNFJZNSnlm3hwYFlzmORfsbmSjWJyBnARg+r7bbIU3Xk20HTpsnbCUA4sEqkzHm1t1x9ZjTp0NmKM5I5HBPtGF4KZE6MXXBLG5yC2a4wIox58UezvGu+RKaFQBKWJdBkqLZ2XQ/tJnpeqNEHfIqUIgH7jnEIZWQemnEzAOSU48pc5oNpqP3QBd/vjf9/85vmrF5/dfvT0zRc//O8/fvXqnz+7efYj+Lj9699uPv6l2qYCowme+kQqNGltb1uNLfCzlISEcU3ZGnuxyyXt9FR2xOkY6YudHLS7vfZBvyXzSUpMxFVJMcWdcRnG2XCsNiEXkkdqJp6gIM4kyU96XODmSpQj7nCXpIjO8+8SVLZ4/fKz1y//fKe0NljWa9k7O4+qdv2trRjchu01XNqK9Qs1QIC3qtajar2eA69g0EgVgNROPc9p2NaOVYcrd5rLvdSRN4FWWXqmCT43V+QUqa3ao9QF8uOEKR1/AiGSslUo5AbaRb/4yZu/PC2gtHU1pTgLKCfbGLGYmQzk1IejiTEeHo/aXWPSHQs/yXg5nJr8a8B69TSeCk4VBDH2v2VAcYdiZlTgRlQxwjODYHdmLCESc3FYcLZNVLXQA6dc1fIgCuOFOpGcQEbdgEyzyUaJDblG8qQ915JOxtoAV7re/P3zTQJ6I52soU3fP/I0VOAoUo6vo/bTf7368sNXL17cPv/5zZc/KMyg10mfBcDPIZryw9TtIeXperNCQSdjacwL67u5+qKJ4iJ1FNKAM8d+LO9QeoRA1Bazyd+oN4eSJyeW9gL6CgW9h1n3muvAdg6RWSSAngsM1TbM7p4pQeXwzK7/+e3vbn/9+e0nz998+Kebj/548/Enr1/+/s2zP6iou336/Panz3SWXy0EUhdxo1VNoCHfR1xDRKMharfx1Qe/MIKQG9ByGLHMSF998GluMqGobE6ymaGlSxUpqrAGzQsLOSOnSqpDQS4VUQ1AW5SwRorQRSxcUDdb5P2qmErEnWR8ozepxowYIXRT34SdFMGZ8H3ltIgqqUfvWHW7rrlKG7GFKWY50x/44RSajIQhHyBWIAWprxfIsHK9g/6w1exnEKXEMPJIJNxQfaCkpRQlpceSUeJqOQpwofFzY1+8Ga3B1lnpzLk0ZurHl7OmJ9Lf5qeKAgK14yhSDVGg3/LH8QKa4eSJ7m6+fLbM9b+HqlfNd8TZuNcp8mGc40IhLLIFQfJlatIslaZ6TLzzqHb2UJgmGwKv8NQJ5tDP9arTuoS+MjJF3ulGURhtTD4ZJ4ENoJOGjGJmFk8x8kxV1+1lZ5UQkoSXfqibn95h2CE+4eXess1MehBePlgWzr6saI8NfU8bs9zVIzVLNkH+QV84yv/7nq+crRlF0FiJB8DSD/DJxXUE992S2qitSEFxA4TV9Qv5Po0YfyIygf5SlNOUcqo61CfiwqU+5PjUaTxSBACY+dnNgppJ6HL1Z5XQ79M5LXkttJL4Lk4CtlwsVAvXK+cporQckmtoMHMzQFKcfg/KhnpHKTObcljIpam8eLtk9HzGyyr2zhQTj0ytqjsl9eqWZ21XdwhpVG0b/sduvW5Zj8TLp54cMgclVyUXMZMDy/7c6fwPwlYzcicdAAA="
$ws.Range("A1").Comment.Text($commentText)

# 3) Update the custom number format (numFmtId 166) applied to the yearly values
$ws.Range("B27:B36").NumberFormat = "###0.000"

# 4) Rename the "Function Description" label to "Function Information"
$ws.Range("A11").Value = "Function Information"

# 5) Tweak the Kurtosis statistic value
$ws.Range("B21").Value = 0.2499825759175085
